$wb = $excel.ActiveWorkbook

$wsSettings  = $wb.Worksheets.Item("Settings")
$wsConstants = $wb.Worksheets.Item("Constants")
$wsAssets    = $wb.Worksheets.Item("Assets")

# --- Settings sheet ---
# B2 "ProcessABCQueue" -> "AutoRecrut"
$wsSettings.Range("B2").Value = "AutoRecrut"
# B3 gets a new value "bello"
$wsSettings.Range("B3").Value = "bello"
$wsSettings.Range("B3").Select()

# --- Constants sheet ---
# New rows 19-22: Business/System exception email subject & body settings
$wsConstants.Range("A19").Value = "BException_Email_Body"
$wsConstants.Range("A20").Value = "BException_Email_Subject "
$wsConstants.Range("A21").Value = "SException_Email_Subject "
$wsConstants.Range("A22").Value = "SException_Email_Body"

$wsConstants.Range("B20").Value = "No Attachment Found"
$wsConstants.Range("B19").Value = "Hello , Kindly note that the Email provided had no attachment added on it ,Kind Regards Admin"
$wsConstants.Range("B21").Value = "Hi , Certain Errors were experienced In the system"
$wsConstants.Range("B22").Value = "System Error"

# Match wrap-text style used by other "description" cells in column B (style index 3 / id "3")
$wsConstants.Range("B19").WrapText = $true
$wsConstants.Range("B20").WrapText = $true

# Trim the last trailing (blank) formatted row at the bottom of the sheet
$wsConstants.Rows.Item(988).Delete()
$wsConstants.Range("A22").Select()

# --- Assets sheet ---
# Trim the last 5 trailing (blank) formatted rows at the bottom of the sheet
$wsAssets.Rows.Item(996).Delete()
$wsAssets.Rows.Item(996).Delete()
$wsAssets.Rows.Item(996).Delete()
$wsAssets.Rows.Item(996).Delete()
$wsAssets.Rows.Item(996).Delete()
$wsAssets.Range("C4").Select()
